# The deck currently uses the "Integral" design theme (slide master -> theme1.xml).
# This commit switches the presentation's theme / colour scheme back to the
# stock "Office Theme" palette (the colours that were still sitting, unused,
# in the deck's second theme part).
#
# PowerPoint's Theme Colors are exposed on the COM object model via
# ThemeColorScheme.Colors(i).RGB (there is no simple "apply theme file"
# call in this host -- custom colour schemes must be edited slot by slot).
# The Colors() index order is the standard OOXML clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# and .RGB is a packed BGR integer (bb*65536 + gg*256 + rr), matching the
# usual VBA/COM RGB() packing.

function ConvertTo-BgrInt($hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($bb * 65536) + ($gg * 256) + $rr
}

# Target "Office Theme" colour scheme (RRGGBB), in clrScheme slot order.
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = ConvertTo-BgrInt $officeColors[$i - 1]
}
